# Replace the 15 lattice-multiplication problems in the 5x3 table with a
# new set of problems/partial-product digits, cell-by-cell, preserving the
# existing run formatting (sz=32) and the <w:br/> line breaks between the
# 5 lines inside each cell ("AA x BB" / "  d    d" / "  ----" / "d|    |" /
# "d|    |"). Word represents a manual line break (<w:br/>) as the
# vertical-tab character (chr 11) inside Range.Text, so building the new
# cell text with that separator and assigning it to Cell.Range.Text
# reproduces the same run/break structure while only touching the text.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$brk = [char]11

$t.Cell(1,1).Range.Text = "66 x 17" + $brk + "  1    7" + $brk + "  ----" + $brk + "6|    |" + $brk + "6|    |"
$t.Cell(1,2).Range.Text = "11 x 94" + $brk + "  9    4" + $brk + "  ----" + $brk + "1|    |" + $brk + "1|    |"
$t.Cell(1,3).Range.Text = "34 x 80" + $brk + "  8    0" + $brk + "  ----" + $brk + "3|    |" + $brk + "4|    |"
$t.Cell(2,1).Range.Text = "38 x 32" + $brk + "  3    2" + $brk + "  ----" + $brk + "3|    |" + $brk + "8|    |"
$t.Cell(2,2).Range.Text = "82 x 24" + $brk + "  2    4" + $brk + "  ----" + $brk + "8|    |" + $brk + "2|    |"
$t.Cell(2,3).Range.Text = "29 x 38" + $brk + "  3    8" + $brk + "  ----" + $brk + "2|    |" + $brk + "9|    |"
$t.Cell(3,1).Range.Text = "20 x 39" + $brk + "  3    9" + $brk + "  ----" + $brk + "2|    |" + $brk + "0|    |"
$t.Cell(3,2).Range.Text = "11 x 71" + $brk + "  7    1" + $brk + "  ----" + $brk + "1|    |" + $brk + "1|    |"
$t.Cell(3,3).Range.Text = "91 x 88" + $brk + "  8    8" + $brk + "  ----" + $brk + "9|    |" + $brk + "1|    |"
$t.Cell(4,1).Range.Text = "54 x 81" + $brk + "  8    1" + $brk + "  ----" + $brk + "5|    |" + $brk + "4|    |"
$t.Cell(4,2).Range.Text = "42 x 71" + $brk + "  7    1" + $brk + "  ----" + $brk + "4|    |" + $brk + "2|    |"
$t.Cell(4,3).Range.Text = "73 x 24" + $brk + "  2    4" + $brk + "  ----" + $brk + "7|    |" + $brk + "3|    |"
$t.Cell(5,1).Range.Text = "35 x 21" + $brk + "  2    1" + $brk + "  ----" + $brk + "3|    |" + $brk + "5|    |"
$t.Cell(5,2).Range.Text = "13 x 28" + $brk + "  2    8" + $brk + "  ----" + $brk + "1|    |" + $brk + "3|    |"
$t.Cell(5,3).Range.Text = "16 x 94" + $brk + "  9    4" + $brk + "  ----" + $brk + "1|    |" + $brk + "6|    |"

Write-Output "Updated 15 lattice multiplication cells"
